# Update the "Förändrad" date column (C) from serial 45203 to 45205
# for every data row (rows 2 through 351) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 351
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)   # Column C
    if ($cell.Value2 -eq 45203) {
        $cell.Value2 = 45205
    }
}
